$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: "T11: 26/3/2020" totals, one more reporting date appended
# after the existing L column ("T10: 25/3/2020").

$ws.Range("M1").Value = "T11: 26/3/2020"

$ws.Range("M2").Value = 1
$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 6
$ws.Range("M5").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 24
$ws.Range("M8").Value = 0
$ws.Range("M9").Value = 31
$ws.Range("M10").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("M14").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 2
$ws.Range("M20").Formula = "=SUM(M2:M19)"

# Match the look of column L (same per-row cell formatting) by copying
# column L's formats onto the new column M.
$ws.Range("L1:L20").Copy()
$ws.Range("M1:M20").PasteSpecial(-4122)

# Widen the new columns like the author did, and leave the selection where
# they left off after typing the new column in.
$ws.Range("L1:L20").ColumnWidth = 15.42
$ws.Range("M1:M20").ColumnWidth = 15.84
[void]$ws.Range("M2").Select()
